$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.346.83'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.441.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +8.05%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.562'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.36%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.503'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.88%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0781'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.25%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.61%  '

$ws.Range("E13").Value = '  +1.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.814.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.442.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.98%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.849'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.08%  '

$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.339.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0941'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.93%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.34%  '

$ws.Range("E24").Value = '  -1.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.03%  '

$ws.Range("E27").Value = '  -1.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -11.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +21.34%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.64'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.74'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0765'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +18.37%  '

$ws.Range("E37").Value = '  -2.31%  '

$ws.Range("E38").Value = '  +0.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.44'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -11.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0295'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.994.60'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '88.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +29.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -12.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.681.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.182'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.43%  '
